# "metaphysics of upper classes and properties"
#
# Expands the workbook from {NamedThing, Person, PersonCollection} into the
# fuller class hierarchy {NamedThing, Observation, NonProcess, MaterialEntity,
# InformationArtifact, Process, Person, PersonCollection}, and adds a
# birth_date column to Person.
#
# NOTE: worksheet variables in this host resolve by tab *position*, not by a
# stable object identity - once the sheet collection is reshuffled (Add /
# Move / Delete), a previously-captured reference can silently repoint to
# whatever now sits at that index. So every insert below re-fetches its
# "insert after" anchor by name right before use, and each new sheet is
# renamed + populated immediately, before any other structural change can
# invalidate the handle.

$wb = $excel.ActiveWorkbook

# --- Observation: raw_value -------------------------------------------------
$anchor = $wb.Worksheets.Item("NamedThing")
$observation = $wb.Worksheets.Add($null, $anchor)
$observation.Name = "Observation"
$observation.Range("A1").Value = "raw_value"

# --- NonProcess: id, name, description --------------------------------------
$anchor = $wb.Worksheets.Item("Observation")
$nonProcess = $wb.Worksheets.Add($null, $anchor)
$nonProcess.Name = "NonProcess"
$nonProcess.Range("A1").Value = "id"
$nonProcess.Range("B1").Value = "name"
$nonProcess.Range("C1").Value = "description"

# --- MaterialEntity: observations, id, name, description --------------------
$anchor = $wb.Worksheets.Item("NonProcess")
$materialEntity = $wb.Worksheets.Add($null, $anchor)
$materialEntity.Name = "MaterialEntity"
$materialEntity.Range("A1").Value = "observations"
$materialEntity.Range("B1").Value = "id"
$materialEntity.Range("C1").Value = "name"
$materialEntity.Range("D1").Value = "description"

# --- InformationArtifact: size_in_bytes, md5, url, id, name, description ----
$anchor = $wb.Worksheets.Item("MaterialEntity")
$informationArtifact = $wb.Worksheets.Add($null, $anchor)
$informationArtifact.Name = "InformationArtifact"
$informationArtifact.Range("A1").Value = "size_in_bytes"
$informationArtifact.Range("B1").Value = "md5"
$informationArtifact.Range("C1").Value = "url"
$informationArtifact.Range("D1").Value = "id"
$informationArtifact.Range("E1").Value = "name"
$informationArtifact.Range("F1").Value = "description"

# --- Process: inputs, outputs, id, name, description -------------------------
$anchor = $wb.Worksheets.Item("InformationArtifact")
$process = $wb.Worksheets.Add($null, $anchor)
$process.Name = "Process"
$process.Range("A1").Value = "inputs"
$process.Range("B1").Value = "outputs"
$process.Range("C1").Value = "id"
$process.Range("D1").Value = "name"
$process.Range("E1").Value = "description"

# --- Person: insert birth_date right after primary_email ---------------------
# (was: primary_email, age_in_years, vital_status, id, name, description)
# The EntireColumn.Insert() shifts the existing vital_status list validation
# from C2:C1048576 to D2:D1048576 automatically.
$person = $wb.Worksheets.Item("Person")
$person.Range("B1").EntireColumn.Insert()
$person.Range("B1").Value = "birth_date"
